$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: correct the "Fallimenti per Fragilità" values for the
# "absolute" locator rows (both LLM and Analitica sections) from 6 to 2.
$ws.Range("E6").Value = 2
$ws.Range("E13").Value = 2

# Recalculate so that dependent formulas (G6, G13, B20, D20, B21, D21) update.
$excel.Calculate()

# Restore the active cell selection to E8, as reflected in the saved view state.
$ws.Range("E8").Select()
